# Apply the cibmtr-reporting-ig update to ValueSet-gvhd-icd10-codes.xlsx
#
# Summary of changes (per the target diff):
#  - Metadata sheet (sheet1): Version 0.1.6 -> 0.1.7, Status active -> draft,
#    Date updated, Publisher text gains a URL, the single "Contact" row is
#    replaced by two "Contact" rows (org + person), and a new "Jurisdiction"
#    row (empty value) is inserted before Description/Purpose/Copyright/
#    Immutable, which shift down by one row.
#  - Include-from sheet (sheet2): content unchanged.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# --- Simple in-place value edits (rows 1-9 keep their row positions) ---
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-27T12:23:18-05:00"
$ws.Range("B9").Value  = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program"

# --- Make room for the new "Jurisdiction" row ---
# Rows 12-15 (Description, Purpose, Copyright, Immutable) need to move down
# to rows 13-16. Extend the existing formatting down to the new row 16 first
# (copy formats only, reusing the existing style instead of creating a new
# one), then rewrite the cell values/labels top-to-bottom.
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "GVHD ICD-10 Codes"

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
